# ReporteDistribucion.xlsx : add "IdEstado"/Distancia column
# A new column is inserted before column P ("CONFIRMACION") on sheet "Informe",
# labeled "Distancia" (defined name DISTANCIA), pushing the existing
# CONFIRMACION / HORARIO / UNREAD_INACTIVE / READ_INACTIVE columns one
# position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informe")

# --- Insert the new column before column P (16) -----------------------
$ws.Columns.Item(16).Insert()

# Column width: mirror the neighbouring (old) column O so the newly
# inserted column gets the same custom width instead of the default one.
$ws.Columns.Item(16).ColumnWidth = $ws.Columns.Item(15).ColumnWidth

# --- Header cell (row 11) / blank data cell (row 12) formatting -------
# The new column must look like the column that used to be at P (now Q),
# so copy its formatting (header style + body style) onto the new cells.
$ws.Range("Q11").Copy()
$ws.Range("P11").PasteSpecial(-4122)

$ws.Range("Q12").Copy()
$ws.Range("P12").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- New header text ----------------------------------------------------
$ws.Range("P11").Value = "Distancia"

# --- Refresh the AutoFilter so it spans the extra column ---------------
$ws.AutoFilterMode = $false
$ws.Range("A11:T11").AutoFilter()

# --- Restore the previous selection/active cell -------------------------
$ws.Range("C8").Select()

# --- Fix up the workbook-level defined names ----------------------------
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Informe!`$A`$11:`$T`$11"
$wb.Names.Item("CONFIRMACION").RefersTo = "=Informe!`$Q`$11"
$wb.Names.Add("DISTANCIA", "=Informe!`$P`$11")
$wb.Names.Item("HORARIO").RefersTo = "=Informe!`$R`$11"
$wb.Names.Item("UNREAD_INACTIVE").RefersTo = "=Informe!`$S`$11"
$wb.Names.Item("READ_INACTIVE").RefersTo = "=Informe!`$T`$11"
